# Apply cryptos list update (price/volume refresh + ranking reshuffle)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '28.250.72'
$ws.Range("E2").Value = '  +2.97%  '

# Row 3
$ws.Range("D3").Value = '1.814.12'
$ws.Range("E3").Value = '  +0.93%  '

# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.26%  '

# Row 5
$ws.Range("D5").Value = '''339.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.48%  '

# Row 6
$ws.Range("D6").Value = '''1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.06%  '

# Row 7
$ws.Range("D7").Value = '''0.3913'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.82%  '

# Row 8
$ws.Range("D8").Value = '''0.3481'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.71%  '

# Row 9
$ws.Range("D9").Value = '''48.39'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.93%  '

# Row 10
$ws.Range("D10").Value = '''1.193'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.99%  '

# Row 11
$ws.Range("D11").Value = '''0.07588'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.88%  '

# Row 12
$ws.Range("D12").Value = '''0.9981'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.36%  '

# Row 13
$ws.Range("D13").Value = '''22.14'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.02%  '

# Row 14
$ws.Range("D14").Value = '''6.514'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.39%  '

# Row 15
$ws.Range("D15").Value = '1.816.47'
$ws.Range("E15").Value = '  +1.06%  '

# Row 16
$ws.Range("D16").Value = '''7.148'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.85%  '

# Row 17
$ws.Range("D17").Value = '''0.00001104'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.27%  '

# Row 18
$ws.Range("D18").Value = '''0.06703'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.61%  '

# Row 19
$ws.Range("D19").Value = '''85.04'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.05%  '

# Row 20
$ws.Range("D20").Value = '''1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.03%  '

# Row 21
$ws.Range("D21").Value = '''17.77'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.98%  '

# Row 22
$ws.Range("D22").Value = '''6.578'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.46%  '

# Row 23
$ws.Range("D23").Value = '28.205.87'
$ws.Range("E23").Value = '  +2.83%  '

# Row 24
$ws.Range("D24").Value = '''12.46'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.87%  '

# Row 25
$ws.Range("D25").Value = '''2.410'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.61%  '

# Row 26
$ws.Range("D26").Value = '''1.497'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.19%  '

# Row 27
$ws.Range("D27").Value = '''2.533'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.37%  '

# Row 28
$ws.Range("D28").Value = '''21.29'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.97%  '

# Row 29
$ws.Range("D29").Value = '''153.89'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.80%  '

# Row 30
$ws.Range("D30").Value = '2.018.10'
$ws.Range("E30").Value = '  +0.84%  '

# Row 31
$ws.Range("D31").Value = '''135.64'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.23%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '''6.159'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.61%  '

# Row 33
$ws.Range("B33").Value = 'HuobiToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D33").Value = '''4.024'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '

# Row 34
$ws.Range("D34").Value = '''0.08854'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.63%  '

# Row 35
$ws.Range("D35").Value = '''13.05'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.94%  '

# Row 36
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").Value = '''5.488'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.40%  '

# Row 37
$ws.Range("B37").Value = 'TheSandbox'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D37").Value = '''0.6952'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.35%  '

# Row 38
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '''0.06558'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.54%  '

# Row 39
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '''0.02425'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.31%  '

# Row 40
$ws.Range("D40").Value = '''1.610'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.57%  '

# Row 41
$ws.Range("D41").Value = '''0.2217'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.63%  '

# Row 42
$ws.Range("D42").Value = '''1.257'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.27%  '

# Row 43
$ws.Range("D43").Value = '''8.508'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.37%  '

# Row 44
$ws.Range("D44").Value = '''14.66'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.83%  '

# Row 45
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '''0.6454'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.27%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").Value = '''3.877'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.09%  '

# Row 47
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").Value = '''2.154'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.74%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '''131.71'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.83%  '

# Row 49
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.07208'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.05%  '

# Row 50
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '''80.08'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.22%  '

# Row 51
$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
$ws.Range("D51").Value = '''1.250'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.45%  '

